$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the terse English values with their "verbose" French equivalents.
$ws.Range("G2").Value = "Restaurant avec cuisine sur place"
$ws.Range("H2").Value = "Concédée"
$ws.Range("I2").Value = "Public"

# G2 (type_production) switches to a Times New Roman font, keeping its
# existing "@" text number format.
$ws.Range("G2").Font.Name = "Times New Roman"
$ws.Range("G2").Font.Size = 10

# H2/I2 (type_gestion / modèle_économique) keep their Arial font but
# switch to General format and wrap their (now longer) text.
$ws.Range("H2:I2").NumberFormat = "General"
$ws.Range("H2:I2").WrapText = $true

# Move the active selection to I2 (matches the saved selection in the sheet).
[void]$ws.Range("I2").Select()
